$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# The "Neighbourhood Road 3" heading paragraph loses its stray trailing
# empty run, and the next paragraph's three runs (split across
# "...regardless I " / "forge on" / ".") collapse into a single run.

# Merge the second paragraph's three runs into one run with the full text.
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start
$p2End = $p2.Range.End
$d.Range($p2Start, $p2End).Text = "Lilith’s place is actually pretty far away, located in a rougher, tucked away part of town. The worn-out buildings combined with the smell of smoke and the grey clouds looming overhead make me feel uneasy, but regardless I forge on."

# Rebuild paragraph 1 ("Neighbourhood Road 3") from scratch so it no longer
# carries the stray trailing empty run: delete the paragraph (its text and
# its paragraph mark) entirely, then insert a clean bold paragraph in its
# place.
$p2 = $d.Paragraphs.Item(2)
$d.Range(0, $p2.Range.Start).Delete()
$d.Range(0, 0).InsertParagraphBefore()
$heading = $d.Paragraphs.Item(1)
$heading.Range.Text = "Neighbourhood Road 3"
$heading.Range.Font.Bold = 1

# --- Edit 2 -----------------------------------------------------------
# Insert a new (italic) "Direction: Screen shakes" paragraph right after
# "Lilith's Dad: Is that so?" (the plain one, no ellipsis) and before
# "He punches me in the stomach...".
$findRange = $d.Content
$found = $findRange.Find.Execute("Lilith’s Dad: Is that so?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
  throw "Could not locate target paragraph for Direction insert"
}
$targetPara = $findRange.Paragraphs.Item(1)
$targetPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetPara.Index + 1)
$dirStart = $newPara.Range.Start
$dirText = "Direction: Screen shakes"
$newPara.Range.Text = $dirText
$d.Range($dirStart, $dirStart + $dirText.Length).Font.Italic = 1

Write-Host "edit complete"
